$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NumeroEntiAttivi (column C) counts and re-rank rows 47-79
# (names/fiscal codes in columns A/B swapped to keep the table sorted
# descending by NumeroEntiAttivi after the count updates).

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = '1083'

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = '608'

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = '511'

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = '429'

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = '420'

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = '362'

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = '252'

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = '249'

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = '227'

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = '203'

$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = '181'

$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = '170'

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = '170'

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = '127'

$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = '122'

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = '82'

$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = '69'

$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = '61'

$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = '51'

$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = '35'

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = '30'

$ws.Range("A47").Value = 'Nexi SpA'

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = '13212880150'

$ws.Range("A48").Value = 'Citta'' Metropolitana di Roma Capitale'

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = '80034390585'

$ws.Range("A51").Value = 'Si.Form Consulting srl'

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = '03943960827'

$ws.Range("A52").Value = 'Servizi Locali SpA'

$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = '03170580751'

$ws.Range("A55").Value = 'UBI Banca'

$ws.Range("B55").NumberFormat = "@"
$ws.Range("B55").Value = '03053920165'

$ws.Range("A56").Value = 'Comune di Catania'

$ws.Range("B56").NumberFormat = "@"
$ws.Range("B56").Value = '00137020871'

$ws.Range("A57").Value = 'ARCA Servizi s.r.l'

$ws.Range("B57").NumberFormat = "@"
$ws.Range("B57").Value = '09106071005'

$ws.Range("C57").NumberFormat = "@"
$ws.Range("C57").Value = '8'

$ws.Range("A58").Value = 'Be Smart s.r.l.'

$ws.Range("B58").NumberFormat = "@"
$ws.Range("B58").Value = '05817461006'

$ws.Range("A60").Value = 'ARGO SOFTWARE SRL'

$ws.Range("B60").NumberFormat = "@"
$ws.Range("B60").Value = '00838520880'

$ws.Range("C60").NumberFormat = "@"
$ws.Range("C60").Value = '5'

$ws.Range("A61").Value = 'Phoenix IT Solutions S.r.L'

$ws.Range("B61").NumberFormat = "@"
$ws.Range("B61").Value = '07623321218'

$ws.Range("A62").Value = 'CityPoste Payment Digital S.r.l.'

$ws.Range("B62").NumberFormat = "@"
$ws.Range("B62").Value = '02003750672'

$ws.Range("A63").Value = 'ISWEB S.p.A.'

$ws.Range("B63").NumberFormat = "@"
$ws.Range("B63").Value = '01722270665'

$ws.Range("A65").Value = 'Linea Comune Spa'

$ws.Range("B65").NumberFormat = "@"
$ws.Range("B65").Value = '05591950489'

$ws.Range("A66").Value = 'KOINE'' SRL'

$ws.Range("B66").NumberFormat = "@"
$ws.Range("B66").Value = '01934790971'

$ws.Range("A67").Value = 'Softline srl'

$ws.Range("B67").NumberFormat = "@"
$ws.Range("B67").Value = '12299030150'

$ws.Range("A69").Value = 'San Marco SPA'

$ws.Range("B69").NumberFormat = "@"
$ws.Range("B69").Value = '04142440728'

$ws.Range("A70").Value = 'BANCA MONTE DEI PASCHI DI SIENA'

$ws.Range("B70").NumberFormat = "@"
$ws.Range("B70").Value = '00884060526'

$ws.Range("A71").Value = 'Società Almaviva S.p.A.'

$ws.Range("B71").NumberFormat = "@"
$ws.Range("B71").Value = '08450891000'

$ws.Range("A72").Value = 'Engineering Ingegneria Informatica SpA'

$ws.Range("B72").NumberFormat = "@"
$ws.Range("B72").Value = '00967720285'

$ws.Range("A73").Value = 'Banco BPM Società per Azioni'

$ws.Range("B73").NumberFormat = "@"
$ws.Range("B73").Value = '09722490969'

$ws.Range("A74").Value = 'Ministero dello Sviluppo Economico'

$ws.Range("B74").NumberFormat = "@"
$ws.Range("B74").Value = '80230390587'

$ws.Range("A75").Value = 'Agenzia Italiana del Farmaco - AIFA'

$ws.Range("B75").NumberFormat = "@"
$ws.Range("B75").Value = '97345810580'

$ws.Range("A76").Value = 'Noviservice srl'

$ws.Range("B76").NumberFormat = "@"
$ws.Range("B76").Value = '02789990922'

$ws.Range("A77").Value = 'Banca Nazionale del Lavoro S.p.A.'

$ws.Range("B77").NumberFormat = "@"
$ws.Range("B77").Value = '09339391006'

$ws.Range("A78").Value = 'MegASP S.r.l.'

$ws.Range("B78").NumberFormat = "@"
$ws.Range("B78").Value = '09898030151'

$ws.Range("A79").Value = 'I.C.A. - Imposte Comunali Affini – s.r.l.'

$ws.Range("B79").NumberFormat = "@"
$ws.Range("B79").Value = '02478610583'
